$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" — copy H1's formatting (bold/bordered
# header style) first, then overwrite the value so the cells pick up the
# existing header style instead of minting a new one.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data rows for the new I/J columns.
$data = @{
    2  = @(5, 6)
    3  = @(5, 5)
    4  = @(6, 6)
    5  = @(6, 6)
    6  = @(5, 5)
    7  = @(1, 2)
    8  = @(8, 8)
    9  = @(3, 4)
    10 = @(5, 5)
    11 = @(4, 4)
    12 = @(2, 3)
    13 = @(2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
